$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15
$ws.Cells.Item(15, 2).Value = 6893654
$ws.Cells.Item(15, 3).Value = 'Germany Oberliga Hamburg'
$ws.Cells.Item(15, 4).Value = (Get-Date -Year 2023 -Month 8 -Day 6 -Hour 9 -Minute 0 -Second 0)
$ws.Cells.Item(15, 5).Value = 'Niendorfer TSV'
$ws.Cells.Item(15, 6).Value = 'Wandsbeker TSV Concordia'
$ws.Cells.Item(15, 7).Value = 6
$ws.Cells.Item(15, 8).Value = 2
$ws.Cells.Item(15, 9).Value = $null
$ws.Cells.Item(15, 10).Value = $null
$ws.Cells.Item(15, 11).Value = 'H'
$ws.Cells.Item(15, 12).Value = 1.444
$ws.Cells.Item(15, 13).Value = 5
$ws.Cells.Item(15, 14).Value = 4.5
$ws.Cells.Item(15, 15).Value = 1.333
$ws.Cells.Item(15, 16).Value = 5.5
$ws.Cells.Item(15, 17).Value = 6
$ws.Cells.Item(15, 18).Value = -1.75
$ws.Cells.Item(15, 19).Value = 1.925
$ws.Cells.Item(15, 20).Value = 1.875
$ws.Cells.Item(15, 21).Value = 4
$ws.Cells.Item(15, 22).Value = 1.875
$ws.Cells.Item(15, 23).Value = 1.925
$ws.Cells.Item(15, 24).Value = 0.333
$ws.Cells.Item(15, 25).Value = -1
$ws.Cells.Item(15, 26).Value = -1
$ws.Cells.Item(15, 27).Value = 0.925
$ws.Cells.Item(15, 28).Value = -1
$ws.Cells.Item(15, 29).Value = 0.875
$ws.Cells.Item(15, 30).Value = -1

# Row 16
$ws.Cells.Item(16, 2).Value = 6893332
$ws.Cells.Item(16, 3).Value = 'Germany Oberliga Hamburg'
$ws.Cells.Item(16, 4).Value = (Get-Date -Year 2023 -Month 8 -Day 6 -Hour 9 -Minute 0 -Second 0)
$ws.Cells.Item(16, 5).Value = 'TSV Buchholz 08'
$ws.Cells.Item(16, 6).Value = 'ETSV Hamburg'
$ws.Cells.Item(16, 7).Value = 4
$ws.Cells.Item(16, 8).Value = 2
$ws.Cells.Item(16, 9).Value = $null
$ws.Cells.Item(16, 10).Value = $null
$ws.Cells.Item(16, 11).Value = 'H'
$ws.Cells.Item(16, 12).Value = 2.3
$ws.Cells.Item(16, 13).Value = 4
$ws.Cells.Item(16, 14).Value = 2.3
$ws.Cells.Item(16, 15).Value = 3.4
$ws.Cells.Item(16, 16).Value = 4.5
$ws.Cells.Item(16, 17).Value = 1.666
$ws.Cells.Item(16, 18).Value = 0.75
$ws.Cells.Item(16, 19).Value = 1.925
$ws.Cells.Item(16, 20).Value = 1.875
$ws.Cells.Item(16, 21).Value = 3.25
$ws.Cells.Item(16, 22).Value = 1.875
$ws.Cells.Item(16, 23).Value = 1.925
$ws.Cells.Item(16, 24).Value = 2.4
$ws.Cells.Item(16, 25).Value = -1
$ws.Cells.Item(16, 26).Value = -1
$ws.Cells.Item(16, 27).Value = 0.925
$ws.Cells.Item(16, 28).Value = -1
$ws.Cells.Item(16, 29).Value = 0.875
$ws.Cells.Item(16, 30).Value = -1

# Row 17
$ws.Cells.Item(17, 2).Value = 6893658
$ws.Cells.Item(17, 3).Value = 'Germany Oberliga Hamburg'
$ws.Cells.Item(17, 4).Value = (Get-Date -Year 2023 -Month 8 -Day 6 -Hour 10 -Minute 0 -Second 0)
$ws.Cells.Item(17, 5).Value = 'FC Union Tornesch'
$ws.Cells.Item(17, 6).Value = 'Tus Dassendorf'
$ws.Cells.Item(17, 7).Value = 1
$ws.Cells.Item(17, 8).Value = 5
$ws.Cells.Item(17, 9).Value = $null
$ws.Cells.Item(17, 10).Value = $null
$ws.Cells.Item(17, 11).Value = 'A'
$ws.Cells.Item(17, 12).Value = 11
$ws.Cells.Item(17, 13).Value = 9
$ws.Cells.Item(17, 14).Value = 1.1
$ws.Cells.Item(17, 15).Value = 19
$ws.Cells.Item(17, 16).Value = 10
$ws.Cells.Item(17, 17).Value = 1.083
$ws.Cells.Item(17, 18).Value = 3
$ws.Cells.Item(17, 19).Value = 2
$ws.Cells.Item(17, 20).Value = 1.8
$ws.Cells.Item(17, 21).Value = 4.5
$ws.Cells.Item(17, 22).Value = 1.875
$ws.Cells.Item(17, 23).Value = 1.925
$ws.Cells.Item(17, 24).Value = -1
$ws.Cells.Item(17, 25).Value = -1
$ws.Cells.Item(17, 26).Value = 0.08299999999999996
$ws.Cells.Item(17, 27).Value = -1
$ws.Cells.Item(17, 28).Value = 0.8
$ws.Cells.Item(17, 29).Value = 0.875
$ws.Cells.Item(17, 30).Value = -1

# Row 18
$ws.Cells.Item(18, 2).Value = 6893657
$ws.Cells.Item(18, 3).Value = 'Germany Oberliga Hamburg'
$ws.Cells.Item(18, 4).Value = (Get-Date -Year 2023 -Month 8 -Day 6 -Hour 10 -Minute 0 -Second 0)
$ws.Cells.Item(18, 5).Value = 'TSV Sasel'
$ws.Cells.Item(18, 6).Value = 'TuRa Harksheide'
$ws.Cells.Item(18, 7).Value = 3
$ws.Cells.Item(18, 8).Value = 1
$ws.Cells.Item(18, 9).Value = 1
$ws.Cells.Item(18, 10).Value = 1
$ws.Cells.Item(18, 11).Value = 'H'
$ws.Cells.Item(18, 12).Value = 1.444
$ws.Cells.Item(18, 13).Value = 4.5
$ws.Cells.Item(18, 14).Value = 5
$ws.Cells.Item(18, 15).Value = 1.25
$ws.Cells.Item(18, 16).Value = 5.5
$ws.Cells.Item(18, 17).Value = 7.5
$ws.Cells.Item(18, 18).Value = -2
$ws.Cells.Item(18, 19).Value = 1.875
$ws.Cells.Item(18, 20).Value = 1.925
$ws.Cells.Item(18, 21).Value = 4.25
$ws.Cells.Item(18, 22).Value = 1.975
$ws.Cells.Item(18, 23).Value = 1.825
$ws.Cells.Item(18, 24).Value = 0.25
$ws.Cells.Item(18, 25).Value = -1
$ws.Cells.Item(18, 26).Value = -1
$ws.Cells.Item(18, 27).Value = 0
$ws.Cells.Item(18, 28).Value = 0
$ws.Cells.Item(18, 29).Value = -0.5
$ws.Cells.Item(18, 30).Value = 0.4125

# Row 25
$ws.Cells.Item(25, 2).Value = 6893664
$ws.Cells.Item(25, 3).Value = 'Germany Oberliga Hamburg'
$ws.Cells.Item(25, 4).Value = (Get-Date -Year 2023 -Month 8 -Day 13 -Hour 9 -Minute 0 -Second 0)
$ws.Cells.Item(25, 5).Value = 'SV Rugenbergen'
$ws.Cells.Item(25, 6).Value = 'FC Trkiye Wilhelmsburg'
$ws.Cells.Item(25, 7).Value = 0
$ws.Cells.Item(25, 8).Value = 2
$ws.Cells.Item(25, 9).Value = 0
$ws.Cells.Item(25, 10).Value = 0
$ws.Cells.Item(25, 11).Value = 'A'
$ws.Cells.Item(25, 12).Value = 2.25
$ws.Cells.Item(25, 13).Value = 4.2
$ws.Cells.Item(25, 14).Value = 2.3
$ws.Cells.Item(25, 15).Value = 1.909
$ws.Cells.Item(25, 16).Value = 4.333
$ws.Cells.Item(25, 17).Value = 2.8
$ws.Cells.Item(25, 18).Value = -0.5
$ws.Cells.Item(25, 19).Value = 1.975
$ws.Cells.Item(25, 20).Value = 1.825
$ws.Cells.Item(25, 21).Value = 3.5
$ws.Cells.Item(25, 22).Value = 1.85
$ws.Cells.Item(25, 23).Value = 1.95
$ws.Cells.Item(25, 24).Value = -1
$ws.Cells.Item(25, 25).Value = -1
$ws.Cells.Item(25, 26).Value = 1.8
$ws.Cells.Item(25, 27).Value = -1
$ws.Cells.Item(25, 28).Value = 0.825
$ws.Cells.Item(25, 29).Value = -1
$ws.Cells.Item(25, 30).Value = 0.95

# Row 26
$ws.Cells.Item(26, 2).Value = 6893663
$ws.Cells.Item(26, 3).Value = 'Germany Oberliga Hamburg'
$ws.Cells.Item(26, 4).Value = (Get-Date -Year 2023 -Month 8 -Day 13 -Hour 9 -Minute 0 -Second 0)
$ws.Cells.Item(26, 5).Value = 'Niendorfer TSV'
$ws.Cells.Item(26, 6).Value = 'FC Union Tornesch'
$ws.Cells.Item(26, 7).Value = 3
$ws.Cells.Item(26, 8).Value = 1
$ws.Cells.Item(26, 9).Value = 1
$ws.Cells.Item(26, 10).Value = 1
$ws.Cells.Item(26, 11).Value = 'H'
$ws.Cells.Item(26, 12).Value = 1.083
$ws.Cells.Item(26, 13).Value = 11
$ws.Cells.Item(26, 14).Value = 13
$ws.Cells.Item(26, 15).Value = 1.083
$ws.Cells.Item(26, 16).Value = 11
$ws.Cells.Item(26, 17).Value = 13
$ws.Cells.Item(26, 18).Value = -3
$ws.Cells.Item(26, 19).Value = 1.9
$ws.Cells.Item(26, 20).Value = 1.9
$ws.Cells.Item(26, 21).Value = 4.5
$ws.Cells.Item(26, 22).Value = 1.925
$ws.Cells.Item(26, 23).Value = 1.775
$ws.Cells.Item(26, 24).Value = 0.08299999999999996
$ws.Cells.Item(26, 25).Value = -1
$ws.Cells.Item(26, 26).Value = -1
$ws.Cells.Item(26, 27).Value = -1
$ws.Cells.Item(26, 28).Value = 0.8999999999999999
$ws.Cells.Item(26, 29).Value = -1
$ws.Cells.Item(26, 30).Value = 0.7749999999999999

# Row 56
$ws.Cells.Item(56, 2).Value = 6893291
$ws.Cells.Item(56, 3).Value = 'Germany Oberliga Hamburg'
$ws.Cells.Item(56, 4).Value = (Get-Date -Year 2023 -Month 9 -Day 8 -Hour 14 -Minute 30 -Second 0)
$ws.Cells.Item(56, 5).Value = 'FC Sderelbe'
$ws.Cells.Item(56, 6).Value = 'TSV Buchholz 08'
$ws.Cells.Item(56, 7).Value = 3
$ws.Cells.Item(56, 8).Value = 3
$ws.Cells.Item(56, 9).Value = 1
$ws.Cells.Item(56, 10).Value = 1
$ws.Cells.Item(56, 11).Value = 'D'
$ws.Cells.Item(56, 12).Value = 1.833
$ws.Cells.Item(56, 13).Value = 4
$ws.Cells.Item(56, 14).Value = 3.1
$ws.Cells.Item(56, 15).Value = 1.363
$ws.Cells.Item(56, 16).Value = 5.25
$ws.Cells.Item(56, 17).Value = 5
$ws.Cells.Item(56, 18).Value = -1.5
$ws.Cells.Item(56, 19).Value = 1.875
$ws.Cells.Item(56, 20).Value = 1.925
$ws.Cells.Item(56, 21).Value = 4.25
$ws.Cells.Item(56, 22).Value = 1.925
$ws.Cells.Item(56, 23).Value = 1.875
$ws.Cells.Item(56, 24).Value = -1
$ws.Cells.Item(56, 25).Value = 4.25
$ws.Cells.Item(56, 26).Value = -1
$ws.Cells.Item(56, 27).Value = -1
$ws.Cells.Item(56, 28).Value = 0.925
$ws.Cells.Item(56, 29).Value = 0.925
$ws.Cells.Item(56, 30).Value = -1

# Row 57
$ws.Cells.Item(57, 2).Value = 6893680
$ws.Cells.Item(57, 3).Value = 'Germany Oberliga Hamburg'
$ws.Cells.Item(57, 4).Value = (Get-Date -Year 2023 -Month 9 -Day 8 -Hour 14 -Minute 30 -Second 0)
$ws.Cells.Item(57, 5).Value = 'SC Viktoria Hamburg'
$ws.Cells.Item(57, 6).Value = 'TSV Sasel'
$ws.Cells.Item(57, 7).Value = 3
$ws.Cells.Item(57, 8).Value = 2
$ws.Cells.Item(57, 9).Value = 2
$ws.Cells.Item(57, 10).Value = 1
$ws.Cells.Item(57, 11).Value = 'H'
$ws.Cells.Item(57, 12).Value = 4.5
$ws.Cells.Item(57, 13).Value = 4.2
$ws.Cells.Item(57, 14).Value = 1.533
$ws.Cells.Item(57, 15).Value = 3.6
$ws.Cells.Item(57, 16).Value = 4
$ws.Cells.Item(57, 17).Value = 1.75
$ws.Cells.Item(57, 18).Value = 0.75
$ws.Cells.Item(57, 19).Value = 1.85
$ws.Cells.Item(57, 20).Value = 1.95
$ws.Cells.Item(57, 21).Value = 3.75
$ws.Cells.Item(57, 22).Value = 1.825
$ws.Cells.Item(57, 23).Value = 1.975
$ws.Cells.Item(57, 24).Value = 2.6
$ws.Cells.Item(57, 25).Value = -1
$ws.Cells.Item(57, 26).Value = -1
$ws.Cells.Item(57, 27).Value = 0.8500000000000001
$ws.Cells.Item(57, 28).Value = -1
$ws.Cells.Item(57, 29).Value = 0.825
$ws.Cells.Item(57, 30).Value = -1

# Row 64
$ws.Cells.Item(64, 2).Value = 6893687
$ws.Cells.Item(64, 3).Value = 'Germany Oberliga Hamburg'
$ws.Cells.Item(64, 4).Value = (Get-Date -Year 2023 -Month 9 -Day 12 -Hour 14 -Minute 30 -Second 0)
$ws.Cells.Item(64, 5).Value = 'FC Union Tornesch'
$ws.Cells.Item(64, 6).Value = 'SV HalstenbekRellingen'
$ws.Cells.Item(64, 7).Value = 1
$ws.Cells.Item(64, 8).Value = 1
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 10).Value = 1
$ws.Cells.Item(64, 11).Value = 'D'
$ws.Cells.Item(64, 12).Value = 3.6
$ws.Cells.Item(64, 13).Value = 4
$ws.Cells.Item(64, 14).Value = 1.7
$ws.Cells.Item(64, 15).Value = 3.6
$ws.Cells.Item(64, 16).Value = 4
$ws.Cells.Item(64, 17).Value = 1.7
$ws.Cells.Item(64, 18).Value = 0.75
$ws.Cells.Item(64, 19).Value = 1.875
$ws.Cells.Item(64, 20).Value = 1.925
$ws.Cells.Item(64, 21).Value = 3.5
$ws.Cells.Item(64, 22).Value = 1.875
$ws.Cells.Item(64, 23).Value = 1.925
$ws.Cells.Item(64, 24).Value = -1
$ws.Cells.Item(64, 25).Value = 3
$ws.Cells.Item(64, 26).Value = -1
$ws.Cells.Item(64, 27).Value = 0.875
$ws.Cells.Item(64, 28).Value = -1
$ws.Cells.Item(64, 29).Value = -1
$ws.Cells.Item(64, 30).Value = 0.925

# Row 65
$ws.Cells.Item(65, 2).Value = 6893685
$ws.Cells.Item(65, 3).Value = 'Germany Oberliga Hamburg'
$ws.Cells.Item(65, 4).Value = (Get-Date -Year 2023 -Month 9 -Day 12 -Hour 14 -Minute 30 -Second 0)
$ws.Cells.Item(65, 5).Value = 'TuRa Harksheide'
$ws.Cells.Item(65, 6).Value = 'ETSV Hamburg'
$ws.Cells.Item(65, 7).Value = 3
$ws.Cells.Item(65, 8).Value = 1
$ws.Cells.Item(65, 9).Value = 2
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 11).Value = 'H'
$ws.Cells.Item(65, 12).Value = 3.4
$ws.Cells.Item(65, 13).Value = 4.5
$ws.Cells.Item(65, 14).Value = 1.666
$ws.Cells.Item(65, 15).Value = 3.4
$ws.Cells.Item(65, 16).Value = 5
$ws.Cells.Item(65, 17).Value = 1.666
$ws.Cells.Item(65, 18).Value = 0.75
$ws.Cells.Item(65, 19).Value = 1.95
$ws.Cells.Item(65, 20).Value = 1.85
$ws.Cells.Item(65, 21).Value = 3.5
$ws.Cells.Item(65, 22).Value = 1.825
$ws.Cells.Item(65, 23).Value = 1.975
$ws.Cells.Item(65, 24).Value = 2.4
$ws.Cells.Item(65, 25).Value = -1
$ws.Cells.Item(65, 26).Value = -1
$ws.Cells.Item(65, 27).Value = 0.95
$ws.Cells.Item(65, 28).Value = -1
$ws.Cells.Item(65, 29).Value = 0.825
$ws.Cells.Item(65, 30).Value = -1

# Row 66
$ws.Cells.Item(66, 2).Value = 6893684
$ws.Cells.Item(66, 3).Value = 'Germany Oberliga Hamburg'
$ws.Cells.Item(66, 4).Value = (Get-Date -Year 2023 -Month 9 -Day 12 -Hour 14 -Minute 30 -Second 0)
$ws.Cells.Item(66, 5).Value = 'FC Alsterbruder'
$ws.Cells.Item(66, 6).Value = 'SC Viktoria Hamburg'
$ws.Cells.Item(66, 7).Value = 7
$ws.Cells.Item(66, 8).Value = 2
$ws.Cells.Item(66, 9).Value = 2
$ws.Cells.Item(66, 10).Value = 0
$ws.Cells.Item(66, 11).Value = 'H'
$ws.Cells.Item(66, 12).Value = 2.7
$ws.Cells.Item(66, 13).Value = 4
$ws.Cells.Item(66, 14).Value = 2
$ws.Cells.Item(66, 15).Value = 3
$ws.Cells.Item(66, 16).Value = 4.2
$ws.Cells.Item(66, 17).Value = 1.95
$ws.Cells.Item(66, 18).Value = 0.5
$ws.Cells.Item(66, 19).Value = 1.825
$ws.Cells.Item(66, 20).Value = 1.975
$ws.Cells.Item(66, 21).Value = 4
$ws.Cells.Item(66, 22).Value = 1.75
$ws.Cells.Item(66, 23).Value = 1.95
$ws.Cells.Item(66, 24).Value = 2
$ws.Cells.Item(66, 25).Value = -1
$ws.Cells.Item(66, 26).Value = -1
$ws.Cells.Item(66, 27).Value = 0.825
$ws.Cells.Item(66, 28).Value = -1
$ws.Cells.Item(66, 29).Value = 0.75
$ws.Cells.Item(66, 30).Value = -1

# Row 67
$ws.Cells.Item(67, 2).Value = 6892996
$ws.Cells.Item(67, 3).Value = 'Germany Oberliga Hamburg'
$ws.Cells.Item(67, 4).Value = (Get-Date -Year 2023 -Month 9 -Day 12 -Hour 14 -Minute 30 -Second 0)
$ws.Cells.Item(67, 5).Value = 'Uhlenhorster SC Paloma'
$ws.Cells.Item(67, 6).Value = 'SV Rugenbergen'
$ws.Cells.Item(67, 7).Value = 2
$ws.Cells.Item(67, 8).Value = 1
$ws.Cells.Item(67, 9).Value = 0
$ws.Cells.Item(67, 10).Value = 0
$ws.Cells.Item(67, 11).Value = 'H'
$ws.Cells.Item(67, 12).Value = 1.45
$ws.Cells.Item(67, 13).Value = 4.75
$ws.Cells.Item(67, 14).Value = 4.5
$ws.Cells.Item(67, 15).Value = 1.444
$ws.Cells.Item(67, 16).Value = 4.75
$ws.Cells.Item(67, 17).Value = 4.75
$ws.Cells.Item(67, 18).Value = -1.25
$ws.Cells.Item(67, 19).Value = 1.925
$ws.Cells.Item(67, 20).Value = 1.875
$ws.Cells.Item(67, 21).Value = 3.75
$ws.Cells.Item(67, 22).Value = 1.95
$ws.Cells.Item(67, 23).Value = 1.85
$ws.Cells.Item(67, 24).Value = 0.444
$ws.Cells.Item(67, 25).Value = -1
$ws.Cells.Item(67, 26).Value = -1
$ws.Cells.Item(67, 27).Value = -0.5
$ws.Cells.Item(67, 28).Value = 0.4375
$ws.Cells.Item(67, 29).Value = -1
$ws.Cells.Item(67, 30).Value = 0.8500000000000001

# Row 101
$ws.Cells.Item(101, 2).Value = 6893296
$ws.Cells.Item(101, 3).Value = 'Germany Oberliga Hamburg'
$ws.Cells.Item(101, 4).Value = (Get-Date -Year 2023 -Month 10 -Day 8 -Hour 10 -Minute 0 -Second 0)
$ws.Cells.Item(101, 5).Value = 'FC Union Tornesch'
$ws.Cells.Item(101, 6).Value = 'TSV Buchholz 08'
$ws.Cells.Item(101, 7).Value = 1
$ws.Cells.Item(101, 8).Value = 0
$ws.Cells.Item(101, 9).Value = 0
$ws.Cells.Item(101, 10).Value = 0
$ws.Cells.Item(101, 11).Value = 'H'
$ws.Cells.Item(101, 12).Value = 4.75
$ws.Cells.Item(101, 13).Value = 4.75
$ws.Cells.Item(101, 14).Value = 1.444
$ws.Cells.Item(101, 15).Value = 5.25
$ws.Cells.Item(101, 16).Value = 5
$ws.Cells.Item(101, 17).Value = 1.4
$ws.Cells.Item(101, 18).Value = 1.5
$ws.Cells.Item(101, 19).Value = 1.875
$ws.Cells.Item(101, 20).Value = 1.925
$ws.Cells.Item(101, 21).Value = 4.25
$ws.Cells.Item(101, 22).Value = 1.975
$ws.Cells.Item(101, 23).Value = 1.825
$ws.Cells.Item(101, 24).Value = 4.25
$ws.Cells.Item(101, 25).Value = -1
$ws.Cells.Item(101, 26).Value = -1
$ws.Cells.Item(101, 27).Value = 0.875
$ws.Cells.Item(101, 28).Value = -1
$ws.Cells.Item(101, 29).Value = -1
$ws.Cells.Item(101, 30).Value = 0.825

# Row 102
$ws.Cells.Item(102, 2).Value = 6893000
$ws.Cells.Item(102, 3).Value = 'Germany Oberliga Hamburg'
$ws.Cells.Item(102, 4).Value = (Get-Date -Year 2023 -Month 10 -Day 8 -Hour 10 -Minute 0 -Second 0)
$ws.Cells.Item(102, 5).Value = 'TSV Sasel'
$ws.Cells.Item(102, 6).Value = 'Uhlenhorster SC Paloma'
$ws.Cells.Item(102, 7).Value = 2
$ws.Cells.Item(102, 8).Value = 1
$ws.Cells.Item(102, 9).Value = 0
$ws.Cells.Item(102, 10).Value = 0
$ws.Cells.Item(102, 11).Value = 'H'
$ws.Cells.Item(102, 12).Value = 1.8
$ws.Cells.Item(102, 13).Value = 4.2
$ws.Cells.Item(102, 14).Value = 3.1
$ws.Cells.Item(102, 15).Value = 2.3
$ws.Cells.Item(102, 16).Value = 4
$ws.Cells.Item(102, 17).Value = 2.3
$ws.Cells.Item(102, 18).Value = 0
$ws.Cells.Item(102, 19).Value = 1.9
$ws.Cells.Item(102, 20).Value = 1.9
$ws.Cells.Item(102, 21).Value = 3.5
$ws.Cells.Item(102, 22).Value = 1.925
$ws.Cells.Item(102, 23).Value = 1.875
$ws.Cells.Item(102, 24).Value = 1.3
$ws.Cells.Item(102, 25).Value = -1
$ws.Cells.Item(102, 26).Value = -1
$ws.Cells.Item(102, 27).Value = 0.8999999999999999
$ws.Cells.Item(102, 28).Value = -1
$ws.Cells.Item(102, 29).Value = -1
$ws.Cells.Item(102, 30).Value = 0.875

# Row 119
$ws.Cells.Item(119, 2).Value = 6893300
$ws.Cells.Item(119, 3).Value = 'Germany Oberliga Hamburg'
$ws.Cells.Item(119, 4).Value = (Get-Date -Year 2023 -Month 11 -Day 3 -Hour 15 -Minute 30 -Second 0)
$ws.Cells.Item(119, 5).Value = 'TuRa Harksheide'
$ws.Cells.Item(119, 6).Value = 'TSV Buchholz 08'
$ws.Cells.Item(119, 7).Value = 3
$ws.Cells.Item(119, 8).Value = 1
$ws.Cells.Item(119, 9).Value = 2
$ws.Cells.Item(119, 10).Value = 1
$ws.Cells.Item(119, 11).Value = 'H'
$ws.Cells.Item(119, 12).Value = 1.909
$ws.Cells.Item(119, 13).Value = 4.333
$ws.Cells.Item(119, 14).Value = 2.8
$ws.Cells.Item(119, 15).Value = 1.909
$ws.Cells.Item(119, 16).Value = 4.2
$ws.Cells.Item(119, 17).Value = 2.875
$ws.Cells.Item(119, 18).Value = -0.5
$ws.Cells.Item(119, 19).Value = 1.975
$ws.Cells.Item(119, 20).Value = 1.825
$ws.Cells.Item(119, 21).Value = 3.25
$ws.Cells.Item(119, 22).Value = 1.925
$ws.Cells.Item(119, 23).Value = 1.875
$ws.Cells.Item(119, 24).Value = 0.909
$ws.Cells.Item(119, 25).Value = -1
$ws.Cells.Item(119, 26).Value = -1
$ws.Cells.Item(119, 27).Value = 0.9750000000000001
$ws.Cells.Item(119, 28).Value = -1
$ws.Cells.Item(119, 29).Value = 0.925
$ws.Cells.Item(119, 30).Value = -1

# Row 120
$ws.Cells.Item(120, 2).Value = 6893727
$ws.Cells.Item(120, 3).Value = 'Germany Oberliga Hamburg'
$ws.Cells.Item(120, 4).Value = (Get-Date -Year 2023 -Month 11 -Day 3 -Hour 15 -Minute 30 -Second 0)
$ws.Cells.Item(120, 5).Value = 'SC Viktoria Hamburg'
$ws.Cells.Item(120, 6).Value = 'Hamburg Eimsbutteler BC'
$ws.Cells.Item(120, 7).Value = 2
$ws.Cells.Item(120, 8).Value = 2
$ws.Cells.Item(120, 9).Value = 1
$ws.Cells.Item(120, 10).Value = 2
$ws.Cells.Item(120, 11).Value = 'D'
$ws.Cells.Item(120, 12).Value = 3.6
$ws.Cells.Item(120, 13).Value = 4.5
$ws.Cells.Item(120, 14).Value = 1.615
$ws.Cells.Item(120, 15).Value = 2.625
$ws.Cells.Item(120, 16).Value = 4.333
$ws.Cells.Item(120, 17).Value = 2.05
$ws.Cells.Item(120, 18).Value = 0.25
$ws.Cells.Item(120, 19).Value = 1.9
$ws.Cells.Item(120, 20).Value = 1.9
$ws.Cells.Item(120, 21).Value = 3.75
$ws.Cells.Item(120, 22).Value = 1.875
$ws.Cells.Item(120, 23).Value = 1.925
$ws.Cells.Item(120, 24).Value = -1
$ws.Cells.Item(120, 25).Value = 3.333
$ws.Cells.Item(120, 26).Value = -1
$ws.Cells.Item(120, 27).Value = 0.45
$ws.Cells.Item(120, 28).Value = -0.5
$ws.Cells.Item(120, 29).Value = 0.4375
$ws.Cells.Item(120, 30).Value = -0.5

# Row 121
$ws.Cells.Item(121, 2).Value = 6893726
$ws.Cells.Item(121, 3).Value = 'Germany Oberliga Hamburg'
$ws.Cells.Item(121, 4).Value = (Get-Date -Year 2023 -Month 11 -Day 3 -Hour 15 -Minute 30 -Second 0)
$ws.Cells.Item(121, 5).Value = 'ETSV Hamburg'
$ws.Cells.Item(121, 6).Value = 'FC Sderelbe'
$ws.Cells.Item(121, 7).Value = 2
$ws.Cells.Item(121, 8).Value = 0
$ws.Cells.Item(121, 9).Value = 0
$ws.Cells.Item(121, 10).Value = 0
$ws.Cells.Item(121, 11).Value = 'H'
$ws.Cells.Item(121, 12).Value = 2.4
$ws.Cells.Item(121, 13).Value = 4
$ws.Cells.Item(121, 14).Value = 2.25
$ws.Cells.Item(121, 15).Value = 1.909
$ws.Cells.Item(121, 16).Value = 3.8
$ws.Cells.Item(121, 17).Value = 3
$ws.Cells.Item(121, 18).Value = -0.5
$ws.Cells.Item(121, 19).Value = 1.95
$ws.Cells.Item(121, 20).Value = 1.85
$ws.Cells.Item(121, 21).Value = 3.75
$ws.Cells.Item(121, 22).Value = 1.75
$ws.Cells.Item(121, 23).Value = 1.95
$ws.Cells.Item(121, 24).Value = 0.909
$ws.Cells.Item(121, 25).Value = -1
$ws.Cells.Item(121, 26).Value = -1
$ws.Cells.Item(121, 27).Value = 0.95
$ws.Cells.Item(121, 28).Value = -1
$ws.Cells.Item(121, 29).Value = -1
$ws.Cells.Item(121, 30).Value = 0.95

# Row 132
$ws.Cells.Item(132, 2).Value = 6893733
$ws.Cells.Item(132, 3).Value = 'Germany Oberliga Hamburg'
$ws.Cells.Item(132, 4).Value = (Get-Date -Year 2023 -Month 11 -Day 12 -Hour 10 -Minute 0 -Second 0)
$ws.Cells.Item(132, 5).Value = 'SV HalstenbekRellingen'
$ws.Cells.Item(132, 6).Value = 'SC Viktoria Hamburg'
$ws.Cells.Item(132, 7).Value = 1
$ws.Cells.Item(132, 8).Value = 2
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 10).Value = 1
$ws.Cells.Item(132, 11).Value = 'A'
$ws.Cells.Item(132, 12).Value = 2.3
$ws.Cells.Item(132, 13).Value = 4.2
$ws.Cells.Item(132, 14).Value = 2.25
$ws.Cells.Item(132, 15).Value = 2.625
$ws.Cells.Item(132, 16).Value = 4.2
$ws.Cells.Item(132, 17).Value = 2
$ws.Cells.Item(132, 18).Value = 0.25
$ws.Cells.Item(132, 19).Value = 1.95
$ws.Cells.Item(132, 20).Value = 1.85
$ws.Cells.Item(132, 21).Value = 3.75
$ws.Cells.Item(132, 22).Value = 1.9
$ws.Cells.Item(132, 23).Value = 1.9
$ws.Cells.Item(132, 24).Value = -1
$ws.Cells.Item(132, 25).Value = -1
$ws.Cells.Item(132, 26).Value = 1
$ws.Cells.Item(132, 27).Value = -1
$ws.Cells.Item(132, 28).Value = 0.8500000000000001
$ws.Cells.Item(132, 29).Value = -1
$ws.Cells.Item(132, 30).Value = 0.8999999999999999

# Row 133
$ws.Cells.Item(133, 2).Value = 6893735
$ws.Cells.Item(133, 3).Value = 'Germany Oberliga Hamburg'
$ws.Cells.Item(133, 4).Value = (Get-Date -Year 2023 -Month 11 -Day 12 -Hour 10 -Minute 0 -Second 0)
$ws.Cells.Item(133, 5).Value = 'FC Union Tornesch'
$ws.Cells.Item(133, 6).Value = 'FC Trkiye Wilhelmsburg'
$ws.Cells.Item(133, 7).Value = 0
$ws.Cells.Item(133, 8).Value = 2
$ws.Cells.Item(133, 9).Value = 0
$ws.Cells.Item(133, 10).Value = 0
$ws.Cells.Item(133, 11).Value = 'A'
$ws.Cells.Item(133, 12).Value = 1.85
$ws.Cells.Item(133, 13).Value = 4
$ws.Cells.Item(133, 14).Value = 3.1
$ws.Cells.Item(133, 15).Value = 1.95
$ws.Cells.Item(133, 16).Value = 3.8
$ws.Cells.Item(133, 17).Value = 2.9
$ws.Cells.Item(133, 18).Value = -0.25
$ws.Cells.Item(133, 19).Value = 1.8
$ws.Cells.Item(133, 20).Value = 2
$ws.Cells.Item(133, 21).Value = 4
$ws.Cells.Item(133, 22).Value = 1.95
$ws.Cells.Item(133, 23).Value = 1.85
$ws.Cells.Item(133, 24).Value = -1
$ws.Cells.Item(133, 25).Value = -1
$ws.Cells.Item(133, 26).Value = 1.9
$ws.Cells.Item(133, 27).Value = -1
$ws.Cells.Item(133, 28).Value = 1
$ws.Cells.Item(133, 29).Value = -1
$ws.Cells.Item(133, 30).Value = 0.8500000000000001

# Row 240
$ws.Cells.Item(240, 2).Value = 6895041
$ws.Cells.Item(240, 3).Value = 'Germany Oberliga Hamburg'
$ws.Cells.Item(240, 4).Value = (Get-Date -Year 2024 -Month 4 -Day 21 -Hour 9 -Minute 0 -Second 0)
$ws.Cells.Item(240, 5).Value = 'SV HalstenbekRellingen'
$ws.Cells.Item(240, 6).Value = 'Niendorfer TSV'
$ws.Cells.Item(240, 7).Value = 1
$ws.Cells.Item(240, 8).Value = 2
$ws.Cells.Item(240, 9).Value = 1
$ws.Cells.Item(240, 10).Value = 1
$ws.Cells.Item(240, 11).Value = 'A'
$ws.Cells.Item(240, 12).Value = 3
$ws.Cells.Item(240, 13).Value = 4.2
$ws.Cells.Item(240, 14).Value = 1.85
$ws.Cells.Item(240, 15).Value = 2.4
$ws.Cells.Item(240, 16).Value = 4.2
$ws.Cells.Item(240, 17).Value = 2.25
$ws.Cells.Item(240, 18).Value = 0
$ws.Cells.Item(240, 19).Value = 1.975
$ws.Cells.Item(240, 20).Value = 1.825
$ws.Cells.Item(240, 21).Value = 4
$ws.Cells.Item(240, 22).Value = 1.75
$ws.Cells.Item(240, 23).Value = 1.95
$ws.Cells.Item(240, 24).Value = -1
$ws.Cells.Item(240, 25).Value = -1
$ws.Cells.Item(240, 26).Value = 1.25
$ws.Cells.Item(240, 27).Value = -1
$ws.Cells.Item(240, 28).Value = 0.825
$ws.Cells.Item(240, 29).Value = -1
$ws.Cells.Item(240, 30).Value = 0.95

# Row 241
$ws.Cells.Item(241, 2).Value = 6895042
$ws.Cells.Item(241, 3).Value = 'Germany Oberliga Hamburg'
$ws.Cells.Item(241, 4).Value = (Get-Date -Year 2024 -Month 4 -Day 21 -Hour 9 -Minute 0 -Second 0)
$ws.Cells.Item(241, 5).Value = 'SV Rugenbergen'
$ws.Cells.Item(241, 6).Value = 'FC Union Tornesch'
$ws.Cells.Item(241, 7).Value = 3
$ws.Cells.Item(241, 8).Value = 0
$ws.Cells.Item(241, 9).Value = 2
$ws.Cells.Item(241, 10).Value = 0
$ws.Cells.Item(241, 11).Value = 'H'
$ws.Cells.Item(241, 12).Value = 1.615
$ws.Cells.Item(241, 13).Value = 4.5
$ws.Cells.Item(241, 14).Value = 3.6
$ws.Cells.Item(241, 15).Value = 1.833
$ws.Cells.Item(241, 16).Value = 4.2
$ws.Cells.Item(241, 17).Value = 3
$ws.Cells.Item(241, 18).Value = -0.5
$ws.Cells.Item(241, 19).Value = 1.9
$ws.Cells.Item(241, 20).Value = 1.9
$ws.Cells.Item(241, 21).Value = 3.5
$ws.Cells.Item(241, 22).Value = 1.825
$ws.Cells.Item(241, 23).Value = 1.975
$ws.Cells.Item(241, 24).Value = 0.833
$ws.Cells.Item(241, 25).Value = -1
$ws.Cells.Item(241, 26).Value = -1
$ws.Cells.Item(241, 27).Value = 0.8999999999999999
$ws.Cells.Item(241, 28).Value = -1
$ws.Cells.Item(241, 29).Value = -1
$ws.Cells.Item(241, 30).Value = 0.9750000000000001

# Row 260
$ws.Cells.Item(260, 2).Value = 6896535
$ws.Cells.Item(260, 3).Value = 'Germany Oberliga Hamburg'
$ws.Cells.Item(260, 4).Value = (Get-Date -Year 2024 -Month 5 -Day 17 -Hour 14 -Minute 0 -Second 0)
$ws.Cells.Item(260, 5).Value = 'Niendorfer TSV'
$ws.Cells.Item(260, 6).Value = 'Uhlenhorster SC Paloma'
$ws.Cells.Item(260, 7).Value = 2
$ws.Cells.Item(260, 8).Value = 2
$ws.Cells.Item(260, 9).Value = 1
$ws.Cells.Item(260, 10).Value = 1
$ws.Cells.Item(260, 11).Value = 'D'
$ws.Cells.Item(260, 12).Value = 1.333
$ws.Cells.Item(260, 13).Value = 5
$ws.Cells.Item(260, 14).Value = 6
$ws.Cells.Item(260, 15).Value = 1.38
$ws.Cells.Item(260, 16).Value = 4.75
$ws.Cells.Item(260, 17).Value = 5.5
$ws.Cells.Item(260, 18).Value = -1.5
$ws.Cells.Item(260, 19).Value = 1.825
$ws.Cells.Item(260, 20).Value = 1.975
$ws.Cells.Item(260, 21).Value = 4.5
$ws.Cells.Item(260, 22).Value = 1.925
$ws.Cells.Item(260, 23).Value = 1.875
$ws.Cells.Item(260, 24).Value = -1
$ws.Cells.Item(260, 25).Value = 3.75
$ws.Cells.Item(260, 26).Value = -1
$ws.Cells.Item(260, 27).Value = -1
$ws.Cells.Item(260, 28).Value = 0.9750000000000001
$ws.Cells.Item(260, 29).Value = -1
$ws.Cells.Item(260, 30).Value = 0.875

# Row 261
$ws.Cells.Item(261, 2).Value = 6895399
$ws.Cells.Item(261, 3).Value = 'Germany Oberliga Hamburg'
$ws.Cells.Item(261, 4).Value = (Get-Date -Year 2024 -Month 5 -Day 17 -Hour 14 -Minute 0 -Second 0)
$ws.Cells.Item(261, 5).Value = 'SV Rugenbergen'
$ws.Cells.Item(261, 6).Value = 'Hamburg Eimsbutteler BC'
$ws.Cells.Item(261, 7).Value = 1
$ws.Cells.Item(261, 8).Value = 1
$ws.Cells.Item(261, 9).Value = 0
$ws.Cells.Item(261, 10).Value = 0
$ws.Cells.Item(261, 11).Value = 'D'
$ws.Cells.Item(261, 12).Value = 2.5
$ws.Cells.Item(261, 13).Value = 4.2
$ws.Cells.Item(261, 14).Value = 2.1
$ws.Cells.Item(261, 15).Value = 3.5
$ws.Cells.Item(261, 16).Value = 4.1
$ws.Cells.Item(261, 17).Value = 1.7
$ws.Cells.Item(261, 18).Value = 0.75
$ws.Cells.Item(261, 19).Value = 1.875
$ws.Cells.Item(261, 20).Value = 1.925
$ws.Cells.Item(261, 21).Value = 3.5
$ws.Cells.Item(261, 22).Value = 1.875
$ws.Cells.Item(261, 23).Value = 1.925
$ws.Cells.Item(261, 24).Value = -1
$ws.Cells.Item(261, 25).Value = 3.1
$ws.Cells.Item(261, 26).Value = -1
$ws.Cells.Item(261, 27).Value = 0.875
$ws.Cells.Item(261, 28).Value = -1
$ws.Cells.Item(261, 29).Value = -1
$ws.Cells.Item(261, 30).Value = 0.925
